$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:K (and the trailing
# blank column) one position to the right, to E:L. This mirrors the
# author's insertion of a new "most-recent fiscal year" column of data
# ahead of the existing years.
$ws.Columns("D").Insert()

# The freshly inserted column D has no formatting; copy the (just
# shifted-right) formatting from column E, which still carries the
# original per-row number formats/styles (date header row, $ numbers,
# etc.), onto the new column D so the new cells inherit the same look.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the new column's on-screen width to its neighbours (all of
# D:K share the same bestFit width in the original workbook).
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the new fiscal-year figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 18313000
$ws.Range("D9").Value2 = 3004000
$ws.Range("D10").Value2 = 15309000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 246000
$ws.Range("D15").Value2 = 802000
$ws.Range("D17").Value2 = 16665000
$ws.Range("D18").Value2 = 1648000
$ws.Range("D20").Value2 = -5000
$ws.Range("D21").Value2 = 2445000
$ws.Range("D22").Value2 = 1004000
$ws.Range("D23").Value2 = 639000
$ws.Range("D24").Value2 = 176000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 463000
$ws.Range("D27").Value2 = 108000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 3000
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 5000
$ws.Range("D33").Value2 = 111000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 111000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 411000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 2785000
$ws.Range("D44").Value2 = 305000
$ws.Range("D45").Value2 = 1135000
$ws.Range("D46").Value2 = 4636000
$ws.Range("D47").Value2 = 1271000
$ws.Range("D48").Value2 = 7148000
$ws.Range("D49").Value2 = 9012000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 342000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 22409000
$ws.Range("D57").Value2 = 1207000
$ws.Range("D58").Value2 = 182000
$ws.Range("D59").Value2 = 2468000
$ws.Range("D60").Value2 = 3857000
$ws.Range("D61").Value2 = 14644000
$ws.Range("D62").Value2 = 1801000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 22528000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -2236000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = -119000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 111000
$ws.Range("D83").Value2 = 802000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 1049000
$ws.Range("D91").Value2 = -617000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -115000
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -1134000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -200000
